$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before the current row 104 (the "note" row), shifting it down to 105.
$ws.Rows.Item(104).Insert()

# Fill in the new data row 104 with the next day's data.
$ws.Cells.Item(104, 1).Value = 43959
$ws.Cells.Item(104, 2).Value = 304
$ws.Cells.Item(104, 3).Value = 35007
$ws.Cells.Item(104, 4).Value = 99
$ws.Cells.Item(104, 5).Value = 7134

# Update the print area to extend by one row (edit the defined name directly so the
# sheet-name reference is serialized unquoted, matching the workbook's existing style).
$printAreaName = $null
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $candidate = $wb.Names.Item($i)
    if ($candidate.Name -like "*Print_Area*") {
        $printAreaName = $candidate
    }
}
if ($printAreaName -eq $null) {
    $printAreaName = $wb.Names.Item(1)
}
$printAreaName.RefersTo = '=相談件数!$A$1:$E$106'

# Move the active selection down to the newly added row (matches author's edit position).
$ws.Range("A104").Select() | Out-Null
